$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "ss-12"
$ws.Range("A16").Value = "ss-13"
$ws.Range("A22").Value = "ss-19"
$ws.Range("A23").Value = "ss-20"
$ws.Range("A24").Value = "ss-21"
$ws.Range("A25").Value = "ss-22"

$ws.Range("H24").Select()
